$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.794.65'
$ws.Range("E2").Value = '  +4.20%  '
$ws.Range("D3").Value = '1.867.22'
$ws.Range("E3").Value = '  +2.86%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '274.25'
$ws.Range("E5").Value = '  -1.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9994'
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5283'
$ws.Range("E7").Value = '  +4.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3394'
$ws.Range("E8").Value = '  -3.91%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06814'
$ws.Range("E9").Value = '  +1.94%  '
$ws.Range("E10").Value = '  +0.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7934'
$ws.Range("E11").Value = '  -3.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07744'
$ws.Range("E12").Value = '  -1.50%  '
$ws.Range("D13").Value = '1.831.31'
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '90.15'
$ws.Range("E14").Value = '  +2.90%  '
$ws.Range("E15").Value = '  +1.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9982'
$ws.Range("E16").Value = '  -0.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.44'
$ws.Range("E17").Value = '  +2.64%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008007'
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9986'
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("D20").Value = '26.832.72'
$ws.Range("E20").Value = '  +4.16%  '
$ws.Range("D21").Value = '2.097.15'
$ws.Range("E21").Value = '  +2.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.710'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.972'
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.107'
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.372'
$ws.Range("E25").Value = '  +5.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '145.61'
$ws.Range("E26").Value = '  +2.25%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.24'
$ws.Range("E27").Value = '  +0.59%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.652'
$ws.Range("E28").Value = '  -0.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '112.59'
$ws.Range("E29").Value = '  +3.15%  '
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.313'
$ws.Range("E31").Value = '  +2.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08872'
$ws.Range("E32").Value = '  +1.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04920'
$ws.Range("E33").Value = '  +1.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.164'
$ws.Range("E34").Value = '  +2.88%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7271'
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.870'
$ws.Range("E36").Value = '  -0.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.217'
$ws.Range("E37").Value = '  +2.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.331'
$ws.Range("E38").Value = '  -1.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01847'
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5093'
$ws.Range("E40").Value = '  -0.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9393'
$ws.Range("E41").Value = '  -3.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '116.39'
$ws.Range("E42").Value = '  +1.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.127'
$ws.Range("E43").Value = '  -1.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.986'
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9982'
$ws.Range("E45").Value = '  -0.29%  '
$ws.Range("E46").Value = '  -2.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1329'
$ws.Range("E47").Value = '  -2.79%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.327'
$ws.Range("E48").Value = '  +1.78%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.10'
$ws.Range("E49").Value = '  -0.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05945'
$ws.Range("E50").Value = '  +2.10%  '
$ws.Range("E51").Value = '  -1.76%  '
